# IPT_inititation_df.xlsx update
#  - "warm_up" sheet: rebase the ipt_init_perc trend to 0 (D113 endpoint -> 0,
#    which ripples through the D101:D112 linear-interpolation formulas)
#  - "policies" sheet: update a handful of ipt_init_perc values to
#    differentiate incidence rates by policy, then re-sort the table by
#    POLICY_ID (column A)
#  - leave the active sheet / selections matching the author's last view

$wb = $excel.ActiveWorkbook

$wsWarmUp   = $wb.Worksheets.Item("warm_up")
$wsPolicies = $wb.Worksheets.Item("policies")

# --- warm_up: zero out the ipt_init_perc trend (D113 is the literal anchor
#     the D101:D112 shared formulas interpolate from) ---
$wsWarmUp.Range("D113").Value = 0

# --- policies: update values to match the new HPV-aligned assumptions
#     (pre-sort positions) ---
$wsPolicies.Range("D3").Value = 0.17    # POLICY_ID=2, G_SET=1, on_art=yes
$wsPolicies.Range("D9").Value = 0.22    # POLICY_ID=2, G_SET=2, on_art=yes
$wsPolicies.Range("D11").Value = 0      # POLICY_ID=1, G_SET=2, on_art=no
$wsPolicies.Range("D12").Value = 0      # POLICY_ID=2, G_SET=2, on_art=no

# --- policies: sort the data table by POLICY_ID ascending ---
$sortRange = $wsPolicies.Range("A2:D13")
$sortKey   = $wsPolicies.Range("A2:A13")
$sortRange.Sort($sortKey)

# --- restore cursor/selection state to match the saved workbook view ---
$wsPolicies.Range("E11").Select()

$wsWarmUp.Activate()
$wsWarmUp.Range("E59").Select()
